# Updates the crypto price table (cols B-E, rows 2-51) to match the
# latest scraped values (GitHub Actions refresh).
#
# Some of the "Price" column values are plain decimal-looking strings
# (e.g. "302.60"). Assigning such a string straight to Range.Value lets
# Excel auto-detect it as a number, which would round-trip it as a
# numeric cell (e.g. 302.60000000000002) instead of the literal text
# that the source data uses. Set-TextValue forces the cell to a Text
# number format before writing the value (preserving the exact string),
# then restores the cell's original style so formatting stays untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

$ws.Range('D2').Value = '43.082.19'
$ws.Range('E2').Value = '  +1.10%  '

$ws.Range('D3').Value = '2.353.52'
$ws.Range('E3').Value = '  +2.47%  '

Set-TextValue $ws.Range('D4') '0.999'
$ws.Range('E4').Value = '  -0.03%  '

Set-TextValue $ws.Range('D5') '302.60'
$ws.Range('E5').Value = '  +0.50%  '

Set-TextValue $ws.Range('D6') '95.66'
$ws.Range('E6').Value = '  +0.04%  '

$ws.Range('E7').Value = '  -0.40%  '

$ws.Range('E8').Value = '  -0.09%  '

Set-TextValue $ws.Range('D9') '0.497'
$ws.Range('E9').Value = '  +0.92%  '

Set-TextValue $ws.Range('D10') '34.07'
$ws.Range('E10').Value = '  -1.38%  '

Set-TextValue $ws.Range('D11') '0.0788'
$ws.Range('E11').Value = '  +0.24%  '

Set-TextValue $ws.Range('D12') '18.67'
$ws.Range('E12').Value = '  -2.62%  '

$ws.Range('E13').Value = '  +3.28%  '

Set-TextValue $ws.Range('D14') '6.72'
$ws.Range('E14').Value = '  -0.39%  '

$ws.Range('D15').Value = '2.720.09'
$ws.Range('E15').Value = '  +2.56%  '

$ws.Range('D16').Value = '2.351.05'
$ws.Range('E16').Value = '  +2.11%  '

Set-TextValue $ws.Range('D17') '0.795'
$ws.Range('E17').Value = '  +1.64%  '

$ws.Range('D18').Value = '43.066.66'
$ws.Range('E18').Value = '  +1.21%  '

Set-TextValue $ws.Range('D19') '12.22'
$ws.Range('E19').Value = '  -0.27%  '

Set-TextValue $ws.Range('D20') '6.27'
$ws.Range('E20').Value = '  +4.54%  '

$ws.Range('E21').Value = '  -0.24%  '

Set-TextValue $ws.Range('D22') '68.07'
$ws.Range('E22').Value = '  +0.48%  '

Set-TextValue $ws.Range('D23') '235.25'
$ws.Range('E23').Value = '  +0.15%  '

Set-TextValue $ws.Range('D24') '2.23'
$ws.Range('E24').Value = '  -1.78%  '

$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D25') '2.43'
$ws.Range('E25').Value = '  +1.51%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D26') '1.00'
$ws.Range('E26').Value = '  -0.13%  '

Set-TextValue $ws.Range('D27') '24.49'
$ws.Range('E27').Value = '  +0.61%  '

$ws.Range('E28').Value = '  -0.06%  '

Set-TextValue $ws.Range('D29') '9.12'
$ws.Range('E29').Value = '  +0.76%  '

Set-TextValue $ws.Range('D30') '31.30'
$ws.Range('E30').Value = '  -2.56%  '

$ws.Range('E31').Value = '  -0.02%  '

Set-TextValue $ws.Range('D32') '5.03'
$ws.Range('E32').Value = '  +1.33%  '

$ws.Range('E33').Value = '  +3.13%  '

Set-TextValue $ws.Range('D34') '17.17'
$ws.Range('E34').Value = '  -1.75%  '

$ws.Range('E35').Value = '  -1.86%  '

$ws.Range('E36').Value = '  +4.79%  '

$ws.Range('E37').Value = '  -1.03%  '

$ws.Range('E38').Value = '  +1.25%  '

Set-TextValue $ws.Range('D39') '22.38'
$ws.Range('E39').Value = '  +10.47%  '

Set-TextValue $ws.Range('D40') '2.76'
$ws.Range('E40').Value = '  +2.22%  '

$ws.Range('E41').Value = '  -0.14%  '

Set-TextValue $ws.Range('D42') '103.84'
$ws.Range('E42').Value = '  -36.91%  '

$ws.Range('D43').Value = '1.943.53'
$ws.Range('E43').Value = '  -1.18%  '

Set-TextValue $ws.Range('D44') '0.0279'
$ws.Range('E44').Value = '  -0.12%  '

$ws.Range('E45').Value = '  +4.22%  '

Set-TextValue $ws.Range('D46') '9.44'
$ws.Range('E46').Value = '  -9.78%  '

Set-TextValue $ws.Range('D47') '2.74'
$ws.Range('E47').Value = '  -0.82%  '

$ws.Range('D48').Value = '2.583.49'
$ws.Range('E48').Value = '  +2.46%  '

Set-TextValue $ws.Range('D49') '52.87'
$ws.Range('E49').Value = '  -0.49%  '

$ws.Range('E50').Value = '  -3.79%  '

Set-TextValue $ws.Range('D51') '72.19'
$ws.Range('E51').Value = '  +1.35%  '
